$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

function Set-SchemeColor($index, $hex) {
  $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
  $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
  $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
  $rgbVal = $r + ($g * 256) + ($b * 65536)
  $cs.Colors($index).RGB = $rgbVal
}

Set-SchemeColor 1 "000000"
Set-SchemeColor 2 "FFFFFF"
Set-SchemeColor 3 "44546A"
Set-SchemeColor 4 "E7E6E6"
Set-SchemeColor 5 "5B9BD5"
Set-SchemeColor 6 "ED7D31"
Set-SchemeColor 7 "A5A5A5"
Set-SchemeColor 8 "FFC000"
Set-SchemeColor 9 "4472C4"
Set-SchemeColor 10 "70AD47"
Set-SchemeColor 11 "0563C1"
Set-SchemeColor 12 "954F72"

Write-Host "done setting colors"
